$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E column values (rows 3-9 and 11-17)
$ws.Range("E3").Value = 144
$ws.Range("E4").Value = 216
$ws.Range("E5").Value = 216
$ws.Range("E6").Value = 504
$ws.Range("E7").Value = 144
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0

$ws.Range("E11").Value = 17280
$ws.Range("E12").Value = 25920
$ws.Range("E13").Value = 25920
$ws.Range("E14").Value = 60480
$ws.Range("E15").Value = 17280
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0

# Update selection to column F
$ws.Columns("F:F").Select()
